$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H ("Label"), pushing the existing "Note" (H->I)
# and "Update Only *" (I->J) columns one slot to the right, carrying their
# values/styles/data-validation/dimension along automatically.
$ws.Columns("H:H").Insert()

# Fill in the new "Label" column's header + row 2 value (row 3 stays blank,
# matching the source row for the other data columns).
$ws.Range("H1").Value = "Label"
$ws.Range("H2").Value = "Updated"

# The inserted column copies row 3's blank-cell formatting from its left
# neighbour (11pt); the other blank cells in that row use the sheet's
# smaller 10pt body font, so line H3 back up with them.
$ws.Range("H3").Font.Size = 10
$ws.Range("H3").Font.Name = "Arial"

# The comment that used to live on the "Update Only *" header (I1) needs to
# move along with that header to its new home at J1.
$commentText = $ws.Range("I1").Comment.Text()
$ws.Range("I1").Comment.Delete()
$ws.Range("J1").AddComment($commentText)

# Restore the original active-cell selection recorded in the target sheet.
$ws.Range("J8").Select()
